$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E (pushes old E,F,G -> F,G,H, inheriting formats)
$ws.Columns.Item(5).Insert() | Out-Null
# Insert a new column at I (after H)
$ws.Columns.Item(9).Insert() | Out-Null

# Header row new values + new data rows 7 & 8
# (ordered to reproduce the original shared-string insertion order)
$ws.Range("I1").Value = "CV 5 fold Train/Val"
$ws.Range("A7").Value = "LGBM"
$ws.Range("B7").Value = "max_depth=10, class_weight='balanced'"
$ws.Range("E1").Value = "Dataset"
$ws.Range("E7").Value = "preprocessed_train_val_Mar13_0130pm_label_enc"
$ws.Range("A8").Value = "CatBoostClassifier"
$ws.Range("B8").Value = "max_depth=10, auto_class_weights='Balanced', n_estimators=10"
$ws.Range("E8").Value = "preprocessed_train_val_Mar13_0130pm_label_enc"

# Give I7:I8 the percent number format used elsewhere in column I (copy format from I6)
$ws.Range("I6").Copy() | Out-Null
$ws.Range("I7:I8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("I7").Value = 0.38839
$ws.Range("I8").Value = 0.3505

# Column widths (best achievable approximations given engine width quantization)
$ws.Columns.Item(4).ColumnWidth = 48.5
$ws.Columns.Item(5).ColumnWidth = 43.83333333333333
$ws.Columns.Item(9).ColumnWidth = 17.0

# Selection on sheet view
$ws.Range("H16").Select() | Out-Null
